# Refresh the cryptocurrency market table (Coin / Link / Price / Volume(1h))
# to match the latest data pulled from coinranking.com, as produced by the
# scheduled "Updated cryptos list" GitHub Actions workflow run.
#
# Rows 2-17 keep the same coin in the same rank position; only Price/Volume
# move. Starting at row 18, WrappedBTC drops out of the tracked top set, so
# every row below shifts up to the next-ranked coin (row 35 / ARBITRUM stays
# pinned in place), and a new coin (RenderToken) is appended at row 51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; D = '26.274.58'; E = '  +0.64%  ' },
    @{ Row = 3; D = '1.664.31'; E = '  +0.56%  ' },
    @{ Row = 4; D = '1.011'; E = '  +0.84%  ' },
    @{ Row = 5; D = '218.81'; E = '  +0.33%  ' },
    @{ Row = 6; E = '  +0.36%  ' },
    @{ Row = 7; D = '1.011'; E = '  +0.79%  ' },
    @{ Row = 8; D = '0.2636'; E = '  +0.88%  ' },
    @{ Row = 9; D = '0.06360'; E = '  +0.35%  ' },
    @{ Row = 10; D = '20.55'; E = '  +0.58%  ' },
    @{ Row = 11; D = '0.07848'; E = '  +0.99%  ' },
    @{ Row = 12; D = '4.564'; E = '  +1.56%  ' },
    @{ Row = 13; D = '1.660.65'; E = '  -0.39%  ' },
    @{ Row = 14; D = '1.893.01'; E = '  +0.57%  ' },
    @{ Row = 15; D = '0.5526'; E = '  +0.97%  ' },
    @{ Row = 16; D = '0.0₅8176'; E = '  +0.16%  ' },
    @{ Row = 17; D = '65.63'; E = '  +0.36%  ' },
    @{ Row = 18; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.010'; E = '  +0.80%  ' },
    @{ Row = 19; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '4.675'; E = '  +2.51%  ' },
    @{ Row = 20; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '192.73'; E = '  -0.05%  ' },
    @{ Row = 21; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '10.22'; E = '  +1.31%  ' },
    @{ Row = 22; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '6.027'; E = '  +0.02%  ' },
    @{ Row = 23; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.012'; E = '  +0.86%  ' },
    @{ Row = 24; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '144.34'; E = '  +1.75%  ' },
    @{ Row = 25; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.1227'; E = '  -1.92%  ' },
    @{ Row = 26; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '7.208'; E = '  -0.85%  ' },
    @{ Row = 27; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '16.06'; E = '  -0.76%  ' },
    @{ Row = 28; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '1.478'; E = '  +2.82%  ' },
    @{ Row = 29; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.05935'; E = '  +0.02%  ' },
    @{ Row = 30; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '1.281'; E = '  -0.07%  ' },
    @{ Row = 31; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '3.590'; E = '  +2.04%  ' },
    @{ Row = 32; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '3.281'; E = '  +1.15%  ' },
    @{ Row = 33; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '1.617'; E = '  +3.08%  ' },
    @{ Row = 34; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '2.826'; E = '  +1.33%  ' },
    @{ Row = 35; D = '0.9597'; E = '  +1.11%  ' },
    @{ Row = 36; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '2.425'; E = '  +0.54%  ' },
    @{ Row = 37; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '0.5800'; E = '  +2.48%  ' },
    @{ Row = 38; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.01604'; E = '  -0.43%  ' },
    @{ Row = 39; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '0.8657'; E = '  +1.97%  ' },
    @{ Row = 40; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '5.865'; E = '  +0.99%  ' },
    @{ Row = 41; B = 'PaxDollar'; C = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D = '1.010'; E = '  +0.80%  ' },
    @{ Row = 42; B = 'Maker'; C = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D = '1.047.00'; E = '  +2.34%  ' },
    @{ Row = 43; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '104.06'; E = '  +1.41%  ' },
    @{ Row = 44; B = 'RocketPoolETH'; C = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'; D = '1.805.72'; E = '  +0.44%  ' },
    @{ Row = 45; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '57.41'; E = '  +0.40%  ' },
    @{ Row = 46; B = 'BabyDogeCoin'; C = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D = '0.0₈106'; E = '  -5.66%  ' },
    @{ Row = 47; B = 'Frax'; C = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D = '1.014'; E = '  +0.67%  ' },
    @{ Row = 48; B = 'Mantle'; C = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D = '0.4383'; E = '  +2.23%  ' },
    @{ Row = 49; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '7.974'; E = '  +1.95%  ' },
    @{ Row = 50; B = 'Cronos'; C = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D = '0.05162'; E = '  +0.19%  ' },
    @{ Row = 51; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '1.435'; E = '  -2.80%  ' }

)

foreach ($r in $rows) {
    if ($r.ContainsKey('B')) {
        $ws.Cells.Item($r.Row, 2).Value = $r.B
    }
    if ($r.ContainsKey('C')) {
        $ws.Cells.Item($r.Row, 3).Value = $r.C
    }
    if ($r.ContainsKey('D')) {
        # Column D ("Price") holds values like "1.010" or "26.274.58" that
        # must stay plain text (matching the source data) rather than being
        # auto-coerced into numbers by Excel's smart-entry parsing.
        $dCell = $ws.Cells.Item($r.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $r.D
    }
    if ($r.ContainsKey('E')) {
        $ws.Cells.Item($r.Row, 5).Value = $r.E
    }
}
